$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").EntireColumn.Insert()
Write-Output "done"
